# ---------------------------------------------------------------------------
# Commit: "Automatic update of files."
# The underlying data source re-ordered survey records 16-22 (record 23 is
# untouched) and appended a brand-new observation as record 24.
# We reproduce this by writing the new per-cell values directly (cell by
# cell, only where the value actually changes) rather than doing a
# generic row-move, since Excel COM would otherwise "helpfully" reformat
# some of the untouched columns.
# ---------------------------------------------------------------------------
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16
$ws.Range("A16").Value = 111814104
$ws.Range("B16").Value = 56398
$ws.Range("E16").Value = 100109
$ws.Range("F16").Value = "Tretåig hackspett"
$ws.Range("G16").Value = "Picoides tridactylus"
$ws.Range("H16").Value = "(Linnaeus, 1758)"
$ws.Range("P16").Value = "åsele 1:1 (åsele 1:1), Ås lm"
$ws.Range("Q16").Value = 610154.5078508666
$ws.Range("R16").Value = 7121460.305022033
$ws.Range("S16").Value = 1
$ws.Range("Z16").Value = "17:23"
$ws.Range("AB16").Value = "17:23"

# Row 17
$ws.Range("A17").Value = 111814688
$ws.Range("B17").Value = 90087
$ws.Range("D17").Value = "LC"
$ws.Range("E17").Value = 3298
$ws.Range("F17").Value = "Trådticka"
$ws.Range("G17").Value = "Climacocystis borealis"
$ws.Range("H17").Value = "(Fr.) Kotl. & Pouzar"
$ws.Range("Q17").Value = 610011.2059644217
$ws.Range("R17").Value = 7121475.688616944
$ws.Range("Z17").Value = "17:55"
$ws.Range("AB17").Value = "17:55"

# Row 18
$ws.Range("A18").Value = 111815114
$ws.Range("B18").Value = 90660
$ws.Range("E18").Value = 4362
$ws.Range("F18").Value = "Blå taggsvamp"
$ws.Range("G18").Value = "Hydnellum caeruleum"
$ws.Range("H18").Value = "(Hornem.) P.Karst."
$ws.Range("P18").Value = "åsele 1:1, Ås lm"
$ws.Range("Q18").Value = 610384.0265214761
$ws.Range("R18").Value = 7121170.261031131
$ws.Range("S18").Value = 5
$ws.Range("Z18").Value = "18:19"
$ws.Range("AB18").Value = "18:19"

# Row 19
$ws.Range("A19").Value = 111815269
$ws.Range("B19").Value = 90666
$ws.Range("D19").Value = "LC"
$ws.Range("E19").Value = 4364
$ws.Range("F19").Value = "Dropptaggsvamp"
$ws.Range("G19").Value = "Hydnellum ferrugineum"
$ws.Range("H19").Value = "(Fr.:Fr.) P. Karst."
$ws.Range("Q19").Value = 610053.7842541422
$ws.Range("R19").Value = 7121273.15248157
$ws.Range("Z19").Value = "18:27"
$ws.Range("AB19").Value = "18:27"

# Row 20
$ws.Range("A20").Value = 111814478
$ws.Range("B20").Value = 77515
$ws.Range("D20").Value = "NT"
$ws.Range("E20").Value = 6425
$ws.Range("F20").Value = "Garnlav"
$ws.Range("G20").Value = "Alectoria sarmentosa"
$ws.Range("H20").Value = "(Ach.) Ach."
$ws.Range("Q20").Value = 610155.3487898401
$ws.Range("R20").Value = 7121461.207019502
$ws.Range("Z20").Value = "17:41"
$ws.Range("AB20").Value = "17:41"

# Row 21
$ws.Range("A21").Value = 111814591
$ws.Range("B21").Value = 77515
$ws.Range("D21").Value = "NT"
$ws.Range("E21").Value = 6425
$ws.Range("F21").Value = "Garnlav"
$ws.Range("G21").Value = "Alectoria sarmentosa"
$ws.Range("H21").Value = "(Ach.) Ach."
$ws.Range("Q21").Value = 610012.4812897337
$ws.Range("R21").Value = 7121464.398116477
$ws.Range("Z21").Value = "17:50"
$ws.Range("AB21").Value = "17:50"

# Row 22
$ws.Range("A22").Value = 111815024
$ws.Range("B22").Value = 56414
$ws.Range("E22").Value = 100049
$ws.Range("F22").Value = "Spillkråka"
$ws.Range("G22").Value = "Dryocopus martius"
$ws.Range("Q22").Value = 609922.1399673244
$ws.Range("R22").Value = 7121488.212810148
$ws.Range("Z22").Value = "18:12"
$ws.Range("AB22").Value = "18:12"

# Row 24 (new row)
$ws.Range("A24").Value = 111881914
$ws.Range("B24").Value = 88489
$ws.Range("C24").Value = "Ovaliderad"
$ws.Range("D24").Value = "NT"
$ws.Range("E24").Value = 1962
$ws.Range("F24").Value = "Vaddporing"
$ws.Range("G24").Value = "Anomoporia kamtschatica"
$ws.Range("H24").Value = "(Parmasto) Bondartseva"
$ws.Range("P24").Value = "Åsele 1:1, Ås lm"
$ws.Range("Q24").Value = 610408.7246031044
$ws.Range("R24").Value = 7121114.457637121
$ws.Range("S24").Value = 25
$ws.Range("T24").Value = "Västerbotten"
$ws.Range("U24").Value = "Åsele"
$ws.Range("V24").Value = "Åsele lappmark"
$ws.Range("W24").Value = "Åsele"
$ws.Range("Z24").Value = "00:00"
$ws.Range("AB24").Value = "00:00"
$ws.Range("AD24").Value = $false
$ws.Range("AE24").Value = $false
$ws.Range("AG24").Value = $false
$ws.Range("AW24").Value = "Ulrika Karlsson"
$ws.Range("AX24").Value = "Ulrika Karlsson"

# Startdatum/Slutdatum ("2023-08-31") must stay literal text, not get
# auto-parsed into an Excel date serial. Copying from an existing cell that
# already holds this exact literal text (General-formatted) sidesteps the
# auto-detect that a plain .Value = "2023-08-31" assignment would trigger.
$ws.Range("Y16").Copy($ws.Range("Y24"))
$ws.Range("Y16").Copy($ws.Range("AA24"))

